# Bugfixed the naive forecaster component module
#
# The "date" column (A2:A82) previously stored Excel serial-date numbers
# (with a custom "YYYY-MM-DD HH:MM:SS" number format) representing the
# 15th of the middle month of each quarter. This replaces those values
# with plain text quarter labels such as "2005Q1", "2005Q2", etc.
# (skipping 2011Q2 and 2014Q2, which were absent from the source series).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the ordered list of quarter labels that belong in A2:A82.
$quarters = New-Object System.Collections.Generic.List[string]
for ($year = 2005; $year -le 2025; $year++) {
    for ($q = 1; $q -le 4; $q++) {
        $label = "$year" + "Q" + "$q"
        if ($label -eq "2011Q2") { continue }
        if ($label -eq "2014Q2") { continue }
        if ($year -eq 2025 -and $q -gt 3) { continue }
        [void]$quarters.Add($label)
    }
}

# Write the quarter-label text into A2:A82 (81 rows).
$row = 2
foreach ($label in $quarters) {
    $ws.Cells.Item($row, 1).Value = $label
    $row++
}

# Re-apply the same (bold, bordered, centered) formatting that the header
# row already uses, so the date column no longer carries the custom
# date-time number format that is being removed from the workbook.
$ws.Range("A1").Copy()
$ws.Range("A2:A82").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
